$d = $word.ActiveDocument

# Locate the paragraph that ends with "...I started to make some UI. " -- this is the
# paragraph that currently holds the trailing _GoBack bookmark and that needs to grow
# two new sibling list-paragraphs after it (with the bookmark ending up in the very
# last of the three).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*started to make some UI*") {
        $targetIndex = $i
    }
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$r = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

# Replace that paragraph's content (in place, using a plain numeric Range so the
# engine merges rather than duplicates it) with: the original two runs, the new
# sentences describing the UI (including the grammar-check proofErr markers around
# the lone "a"), then two brand-new ListParagraph-styled paragraphs -- the second
# of the two script/score sentences, and a final empty paragraph that now owns the
# relocated _GoBack bookmark.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="003A6799" w:rsidRDefault="003A6799" w:rsidP="00227E11">
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>A</w:t></w:r>
<w:r><w:t xml:space="preserve">fter that I started to make some UI. </w:t></w:r>
<w:r><w:t xml:space="preserve">Including a scene for when first start the game which contain a start game button and </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>a</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> exit button</w:t></w:r>
<w:r><w:t>, a in-game menu when the player dies that include a play again button and a return button. A score UI and meters travelled UI is also included</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr>
<w:r><w:t xml:space="preserve">Then I included in the script </w:t></w:r>
<w:r><w:t>to store the score and the meters.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr>
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
